$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Bo Mang Nha" gallery link added next to the header row
$ws.Range("G2").Value = "https://drive.google.com/drive/folders/16sW6i0hXS5sGh2gQ7PkoB9GitKugGoNq?usp=sharing"

# Rename / shorten the 2026 "Bon Mang Nha" + "Tu Lieu" activity texts
$ws.Range("E7").Value = "Bổn Mạng Nhà"
$ws.Range("E8").Value = "Tư Liệu"
$ws.Range("F8").Value = "gửi ảnh vào link "

# Turn the G8 photo-upload link into an actual clickable hyperlink
$ws.Hyperlinks.Add($ws.Range("G8"), "https://drive.google.com/drive/folders/1IVawCMt9xO_6Cnvzh2S28Q6U66pePz7e?usp=sharing")
$ws.Range("G8").Style = "Hyperlink"

# Restore the last-used selection recorded in the sheet view
$ws.Range("L14").Select() | Out-Null
